# BOM updates - add IF synth, filters, mixer
#
# - Rename sheet1 (drop the "-v1" suffix)
# - Break/remove the external workbook link (synthesizers.xlsx) so the
#   <externalReferences>/externalLink1.xml parts are dropped on save
# - Update the "Source: ...sch" note on sheet1 to the new generic name
# - Update the revision-history note on sheet1 to call out boards 34-40
# - Move the sheet1 selection from B21 -> A20

$wb = $excel.ActiveWorkbook

# --- sheet1: rename + content tweaks -------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "microstrip-test-4-layer-sma"

$ws1.Range("A2").Value = "Source: microstrip-test-4-layer-sma-xx-xx.sch"
$ws1.Range("C18").Value = "Initial BOM. Assembled 2x boards for 34-40"

[void]$ws1.Range("A20").Select()

# --- drop the external reference to synthesizers.xlsx ---------------------
foreach ($src in @($wb.LinkSources())) {
    [void]$wb.BreakLink($src, 1)
}
